$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing BOM line-item quantities ---
$ws.Range("G4").Value2 = 19
$ws.Range("G21").Value2 = 3
$ws.Range("G24").Value2 = 11
$ws.Range("G25").Value2 = 3
$ws.Range("G35").Value2 = 1

# --- Add new BOM row 47: "RTC backup cap" (Seiko CPH3225A) ---
# Copy formatting from row 46 (the last existing data row) so the new
# row matches the established table styling (fonts/borders/fills).
$ws.Range("A46:B46").Copy()
$ws.Range("A47:B47").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("H46:I46").Copy()
$ws.Range("H47:I47").PasteSpecial(-4122)   # xlPasteFormats

# The Quantity column (G) on row 46 has no border, so mirror that for G47.
$ws.Range("A46").Copy()
$ws.Range("G47").PasteSpecial(-4122)       # xlPasteFormats
$ws.Range("G47").Borders.LineStyle = -4142 # xlLineStyleNone

$excel.CutCopyMode = $false

$ws.Range("A47").Value2 = "RTC backup cap"
$ws.Range("B47").Value2 = "11mF"
$ws.Range("G47").Value2 = 1
$ws.Range("H47").Value2 = "CPH3225A"
$ws.Range("I47").Formula = '=HYPERLINK("http://www.digikey.ca/product-detail/en/seiko-instruments/CPH3225A/728-1067-1-ND/4747400","Digikey - 728-1067-1-ND")'

# --- Update the saved view state (scroll position / zoom / selection) ---
$win = $excel.ActiveWindow
$win.Zoom = 100
$ws.Range("G36").Select() | Out-Null
